$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.535.03"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.53%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.826.04"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.45%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.65"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.33%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5192"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -4.54%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3936"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +3.87%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07715"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +3.28%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.00"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.22%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.86%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.04"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +3.06%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.002"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.558"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.51%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.824.39"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.56%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.25"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +4.43%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06616"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.43%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.61%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.002"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.057"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.542.79"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.46%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.28%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "156.88"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.30%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.035.54"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.44%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.423"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +3.73%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.10"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.70%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.135"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.98%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1107"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.658"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.14%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.54%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07230"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +5.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2247"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.74%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.995"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +6.63%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02334"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.163"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.91%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6253"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.67%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.188"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.40%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.15%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.70%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.40"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5910"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +3.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.706"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.56%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.61"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.10%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.984"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +3.60%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.187"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.40%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06940"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.85%  "
